$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# SECTION 1 (altaProveedor table)
# ---------------------------------------------------------------------------

# 1a. Shorten the "apellido" description text, dropping the "correo"/"descripcion"
#     tails which will move into their own paragraphs.
$null = $d.Content.Find.Execute(
    ": IDEM, pero corresponde al apellido. correo: IDEM anterior, pero corresponde a la dirección de correo electrónico. descripcion: IDEM anterior, pero corresponde a una descripción general",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": IDEM, pero corresponde al apellido. ", 2)

# 1b. Insert a brand-new paragraph for "correo" right after the "apellido" one,
#     before the existing "descripción" paragraph. (Identified by the paragraph
#     that immediately follows starting with "descripci\u00f3n:".)
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "apellido: IDEM, pero corresponde al apellido. `r" -and $i -lt $count) {
        $next = $d.Paragraphs($i + 1)
        if ($next.Range.Text.StartsWith("descripci")) {
            $insPoint = $d.Range($p.Range.End, $p.Range.End)
            $xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Standard"/>
    <w:autoSpaceDE/>
    <w:spacing w:after="200" w:line="276" w:lineRule="auto"/>
    <w:jc w:val="left"/>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
      <w:color w:val="ED7D31" w:themeColor="accent2"/>
    </w:rPr>
    <w:t>correo</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
    </w:rPr>
    <w:t xml:space="preserve">: IDEM, pero corresponde a la dirección de correo electrónico. descripcion: </w:t>
  </w:r>
</w:p>
'@
            $insPoint.InsertXML($xml)
            break
        }
    }
}

# 1c. The old "descripción" paragraph used to finish with the "url" explanation;
#     now it should just end with a single trailing space. Rewrite the whole
#     paragraph (over its full Range, including the paragraph mark) so that the
#     pre-existing run boundaries are preserved instead of being coalesced.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.StartsWith("descripci") -and $t -match "sitio web") {
        $xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Standard"/>
    <w:autoSpaceDE/>
    <w:spacing w:after="200" w:line="276" w:lineRule="auto"/>
    <w:jc w:val="left"/>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
      <w:color w:val="ED7D31" w:themeColor="accent2"/>
    </w:rPr>
    <w:t>d</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
      <w:color w:val="ED7D31" w:themeColor="accent2"/>
    </w:rPr>
    <w:t>escripción</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
    </w:rPr>
    <w:t xml:space="preserve">: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
    </w:rPr>
    <w:t>IDEM, pero corresponde a</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
    </w:rPr>
    <w:t xml:space="preserve"> una descripción general</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
'@
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# SECTION 2 (altaTurista table)
# ---------------------------------------------------------------------------

# 2a. "nombre" description: "IDEM anterior" -> "IDEM"
$null = $d.Content.Find.Execute(
    ": IDEM anterior, pero corresponde al nombre. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": IDEM, pero corresponde al nombre. ", 2)

# 2b. Shorten "apellido" description, dropping the "correo" tail.
$null = $d.Content.Find.Execute(
    ": IDEM anterior, pero corresponde al apellido. correo: IDEM anterior, pero corresponde a la dirección de correo electrónico.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": IDEM, pero corresponde al apellido. ", 2)

# 2c. Insert a brand-new paragraph for "correo" right after the "apellido" one,
#     before the existing "Nacionalidad" paragraph. (Identified by the paragraph
#     that immediately follows starting with "Nacionalidad".)
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "apellido: IDEM, pero corresponde al apellido. `r" -and $i -lt $count) {
        $next = $d.Paragraphs($i + 1)
        if ($next.Range.Text.StartsWith("Nacionalidad")) {
            $insPoint = $d.Range($p.Range.End, $p.Range.End)
            $xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Standard"/>
    <w:autoSpaceDE/>
    <w:spacing w:after="200" w:line="276" w:lineRule="auto"/>
    <w:jc w:val="left"/>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
      <w:color w:val="ED7D31" w:themeColor="accent2"/>
    </w:rPr>
    <w:t>correo</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Courier New"/>
    </w:rPr>
    <w:t>: IDEM, pero corresponde a la dirección de correo electrónico.</w:t>
  </w:r>
</w:p>
'@
            $insPoint.InsertXML($xml)
            break
        }
    }
}

# 2d. "Nacionalidad" description: "IDEM anterior" -> "IDEM"
$null = $d.Content.Find.Execute(
    ": IDEM anterior, pero corresponde a la nacionalidad.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ": IDEM, pero corresponde a la nacionalidad.", 2)
